# RPA datasets push 2023-12-12
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

# Insert a new data row at row 4 (pushes existing rows 4-21 down to 5-22)
$ws.Rows.Item(4).Insert()

# New IPO entry: 현대힘스
$ws.Range("A4").Value = "현대힘스"
$ws.Range("B4").Value = "2024.01.08~01.12"
$ws.Range("C4").Value = "5,000~6,300"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = 43535
$ws.Range("F4").Value = "미래에셋증권"

# Update the 확정공모가 for IBKS스팩23호 (now row 8) from "-" to the text "2000"
# (use a text-formula + paste-values round-trip so the cell stays plain text
# without leaving a custom number-format behind)
$ws.Range("D8").Formula = "=""2000"""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)

# Drop the last row (스톰테크), now row 22 after the insert above
$ws.Rows.Item(22).Delete()
